$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column W = pib, Column X = quartile
$ws.Range("W2").Value = 30861.45
$ws.Range("W3").Value = 87608.38
$ws.Range("W4").Value = 309613.32
$ws.Range("W5").Value = 214081.99
$ws.Range("W6").Value = 47076.8
$ws.Range("W7").Value = 28885.1
$ws.Range("W8").Value = 188411.03
$ws.Range("X8").Value = "third quartile"
$ws.Range("W9").Value = 120517.63
$ws.Range("W10").Value = 114987.25
$ws.Range("W11").Value = 1934032.42
$ws.Range("W12").Value = 634454.0600000001
$ws.Range("W13").Value = 12801.68
$ws.Range("W14").Value = 50041.66
$ws.Range("W15").Value = 132225.73
$ws.Range("W16").Value = 42710.61
$ws.Range("W17").Value = 1604021.07
$ws.Range("W18").Value = 40219.7
$ws.Range("W19").Value = 48927.17
$ws.Range("W20").Value = 58949.94
$ws.Range("W21").Value = 47916.46
$ws.Range("W22").Value = 189664.68
$ws.Range("X22").Value = "more third quartile"
$ws.Range("W23").Value = 108451.44
